# Update the cryptocurrency price/volume table with the latest scrape.
# A new coin (WrappedeETH) was inserted at row 28, shifting all subsequent
# rows down by one; dogwifhat (previously the last row) drops off the bottom.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force the cell to Text format before assigning, so numeric-looking
# strings (e.g. "93.396.85", "1.00", "0.0000244") are not reinterpreted as
# numbers/dates/scientific notation by Excel.
function Set-TextCell {
    param($Sheet, $Addr, $Val)
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
}

# Row 2
Set-TextCell $ws 'D2' '93.396.85'
Set-TextCell $ws 'E2' '  -5.35%  '
# Row 3
Set-TextCell $ws 'D3' '3.366.69'
Set-TextCell $ws 'E3' '  -3.03%  '
# Row 4
Set-TextCell $ws 'E4' '  +0.10%  '
# Row 5
Set-TextCell $ws 'D5' '234.44'
Set-TextCell $ws 'E5' '  -8.34%  '
# Row 6
Set-TextCell $ws 'D6' '627.71'
Set-TextCell $ws 'E6' '  -6.77%  '
# Row 7
Set-TextCell $ws 'D7' '1.39'
Set-TextCell $ws 'E7' '  -7.80%  '
# Row 8
Set-TextCell $ws 'D8' '0.389'
Set-TextCell $ws 'E8' '  -10.59%  '
# Row 9
Set-TextCell $ws 'E9' '  +0.12%  '
# Row 10
Set-TextCell $ws 'D10' '0.942'
Set-TextCell $ws 'E10' '  -11.50%  '
# Row 11
Set-TextCell $ws 'D11' '3.364.98'
Set-TextCell $ws 'E11' '  -3.02%  '
# Row 12
Set-TextCell $ws 'D12' '0.195'
Set-TextCell $ws 'E12' '  -7.78%  '
# Row 13
Set-TextCell $ws 'D13' '40.39'
Set-TextCell $ws 'E13' '  -13.09%  '
# Row 14
Set-TextCell $ws 'D14' '6.01'
Set-TextCell $ws 'E14' '  -3.59%  '
# Row 15
Set-TextCell $ws 'D15' '93.231.24'
Set-TextCell $ws 'E15' '  -5.35%  '
# Row 16
Set-TextCell $ws 'D16' '3.991.69'
Set-TextCell $ws 'E16' '  -3.42%  '
# Row 17
Set-TextCell $ws 'D17' '0.0000244'
Set-TextCell $ws 'E17' '  -6.63%  '
# Row 18
Set-TextCell $ws 'D18' '8.00'
Set-TextCell $ws 'E18' '  -12.31%  '
# Row 19
Set-TextCell $ws 'D19' '3.372.59'
Set-TextCell $ws 'E19' '  -2.74%  '
# Row 20
Set-TextCell $ws 'D20' '16.89'
Set-TextCell $ws 'E20' '  -8.98%  '
# Row 21
Set-TextCell $ws 'D21' '10.91'
Set-TextCell $ws 'E21' '  -7.50%  '
# Row 22
Set-TextCell $ws 'D22' '491.31'
Set-TextCell $ws 'E22' '  -6.32%  '
# Row 23
Set-TextCell $ws 'D23' '0.451'
Set-TextCell $ws 'E23' '  -16.89%  '
# Row 24
Set-TextCell $ws 'D24' '3.14'
Set-TextCell $ws 'E24' '  -9.48%  '
# Row 25
Set-TextCell $ws 'D25' '0.0000186'
Set-TextCell $ws 'E25' '  -9.39%  '
# Row 26
Set-TextCell $ws 'D26' '6.31'
Set-TextCell $ws 'E26' '  -8.60%  '
# Row 27
Set-TextCell $ws 'D27' '89.96'
Set-TextCell $ws 'E27' '  -8.63%  '
# Row 28
Set-TextCell $ws 'B28' 'WrappedeETH'
Set-TextCell $ws 'C28' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextCell $ws 'D28' '3.551.18'
Set-TextCell $ws 'E28' '  -2.84%  '
# Row 29
Set-TextCell $ws 'B29' 'Aptos'
Set-TextCell $ws 'C29' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D29' '11.50'
Set-TextCell $ws 'E29' '  -9.77%  '
# Row 30
Set-TextCell $ws 'B30' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C30' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D30' '11.34'
Set-TextCell $ws 'E30' '  -8.14%  '
# Row 31
Set-TextCell $ws 'B31' 'Dai'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D31' '1.00'
Set-TextCell $ws 'E31' '  +0.08%  '
# Row 32
Set-TextCell $ws 'B32' 'PancakeSwap'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell $ws 'D32' '2.65'
Set-TextCell $ws 'E32' '  -13.52%  '
# Row 33
Set-TextCell $ws 'B33' 'Hedera'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws 'D33' '0.131'
Set-TextCell $ws 'E33' '  -10.48%  '
# Row 34
Set-TextCell $ws 'B34' 'Binance-PegBSC-USD'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell $ws 'D34' '1.00'
Set-TextCell $ws 'E34' '  -0.16%  '
# Row 35
Set-TextCell $ws 'B35' 'Cronos'
Set-TextCell $ws 'C35' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell $ws 'D35' '0.172'
Set-TextCell $ws 'E35' '  -11.00%  '
# Row 36
Set-TextCell $ws 'B36' 'EthereumClassic'
Set-TextCell $ws 'C36' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell $ws 'D36' '28.66'
Set-TextCell $ws 'E36' '  -5.13%  '
# Row 37
Set-TextCell $ws 'B37' 'PolygonEcosystemToken'
Set-TextCell $ws 'C37' 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextCell $ws 'D37' '0.530'
Set-TextCell $ws 'E37' '  -8.69%  '
# Row 38
Set-TextCell $ws 'B38' 'RenderToken'
Set-TextCell $ws 'C38' 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextCell $ws 'D38' '7.49'
Set-TextCell $ws 'E38' '  -8.68%  '
# Row 39
Set-TextCell $ws 'B39' 'Bittensor'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell $ws 'D39' '523.25'
Set-TextCell $ws 'E39' '  -2.64%  '
# Row 40
Set-TextCell $ws 'E40' '  +0.00%  '
# Row 41
Set-TextCell $ws 'B41' 'Fetch.AI'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell $ws 'D41' '1.39'
Set-TextCell $ws 'E41' '  -9.15%  '
# Row 42
Set-TextCell $ws 'B42' 'Kaspa'
Set-TextCell $ws 'C42' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws 'D42' '0.148'
Set-TextCell $ws 'E42' '  -5.65%  '
# Row 43
Set-TextCell $ws 'B43' 'ARBITRUM'
Set-TextCell $ws 'C43' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws 'D43' '0.880'
Set-TextCell $ws 'E43' '  -1.18%  '
# Row 44
Set-TextCell $ws 'B44' 'WhiteBITCoin'
Set-TextCell $ws 'C44' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextCell $ws 'D44' '24.02'
Set-TextCell $ws 'E44' '  -1.74%  '
# Row 45
Set-TextCell $ws 'B45' 'MantraDAO'
Set-TextCell $ws 'C45' 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextCell $ws 'D45' '3.62'
Set-TextCell $ws 'E45' '  -1.58%  '
# Row 46
Set-TextCell $ws 'B46' 'ImmutableX'
Set-TextCell $ws 'C46' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D46' '1.66'
Set-TextCell $ws 'E46' '  -7.34%  '
# Row 47
Set-TextCell $ws 'B47' 'Filecoin'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D47' '5.50'
Set-TextCell $ws 'E47' '  -5.27%  '
# Row 48
Set-TextCell $ws 'B48' 'Stacks'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D48' '2.15'
Set-TextCell $ws 'E48' '  -4.49%  '
# Row 49
Set-TextCell $ws 'B49' 'VeChain'
Set-TextCell $ws 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws 'D49' '0.0394'
Set-TextCell $ws 'E49' '  -10.72%  '
# Row 50
Set-TextCell $ws 'B50' 'OKB'
Set-TextCell $ws 'C50' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell $ws 'D50' '52.42'
Set-TextCell $ws 'E50' '  -6.41%  '
# Row 51
Set-TextCell $ws 'B51' 'Cosmos'
Set-TextCell $ws 'C51' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws 'D51' '7.96'
Set-TextCell $ws 'E51' '  -9.52%  '
